# Auto-generated edit script: update market-price / profit columns (H-N)
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 88.5625  # ALC!H5 (was 98)
$ws.Cells.Item(5, 9).Value = 88  # ALC!I5 (was 98)
$ws.Cells.Item(5, 10).Value = 92.5  # ALC!J5 (was 0)
$ws.Cells.Item(5, 11).Value = 88  # ALC!K5 (was 98)
$ws.Cells.Item(5, 12).Value = 92.5  # ALC!L5 (was 0)
$ws.Cells.Item(5, 13).Value = 27  # ALC!M5 (was 17)
$ws.Cells.Item(5, 14).Value = -322.5  # ALC!N5 (was None)

$ws.Cells.Item(64, 8).Value = 4552.6313  # ALC!H64 (was 4340.909)
$ws.Cells.Item(64, 9).Value = 4500  # ALC!I64 (was 4285.7144)
$ws.Cells.Item(64, 11).Value = 4500  # ALC!K64 (was 4285.7144)
$ws.Cells.Item(64, 13).Value = -4252  # ALC!M64 (was -4037.7144)

$ws.Cells.Item(67, 8).Value = 4552.6313  # ALC!H67 (was 4340.909)
$ws.Cells.Item(67, 9).Value = 4500  # ALC!I67 (was 4285.7144)
$ws.Cells.Item(67, 11).Value = 4500  # ALC!K67 (was 4285.7144)
$ws.Cells.Item(67, 13).Value = -3642  # ALC!M67 (was -3427.7144)

$ws.Cells.Item(70, 8).Value = 1967.8  # ALC!H70 (was 3023.75)
$ws.Cells.Item(70, 9).Value = 896  # ALC!I70 (was 850)
$ws.Cells.Item(70, 10).Value = 2235.75  # ALC!J70 (was 3748.3333)
$ws.Cells.Item(70, 11).Value = 2688  # ALC!K70 (was 2550)
$ws.Cells.Item(70, 12).Value = 6707.25  # ALC!L70 (was 11244.9999)
$ws.Cells.Item(70, 13).Value = -2418  # ALC!M70 (was -2280)
$ws.Cells.Item(70, 14).Value = -7247.25  # ALC!N70 (was -11784.9999)

$ws.Cells.Item(73, 8).Value = 1967.8  # ALC!H73 (was 3023.75)
$ws.Cells.Item(73, 9).Value = 896  # ALC!I73 (was 850)
$ws.Cells.Item(73, 10).Value = 2235.75  # ALC!J73 (was 3748.3333)
$ws.Cells.Item(73, 11).Value = 2688  # ALC!K73 (was 2550)
$ws.Cells.Item(73, 12).Value = 6707.25  # ALC!L73 (was 11244.9999)
$ws.Cells.Item(73, 13).Value = -1752  # ALC!M73 (was -1614)
$ws.Cells.Item(73, 14).Value = -8579.25  # ALC!N73 (was -13116.9999)

$ws.Cells.Item(76, 8).Value = 55558364  # ALC!H76 (was 57695116)
$ws.Cells.Item(76, 9).Value = 65220120  # ALC!I76 (was 68184540)
$ws.Cells.Item(76, 11).Value = 65220120  # ALC!K76 (was 68184540)
$ws.Cells.Item(76, 13).Value = -65219805  # ALC!M76 (was -68184225)

$ws.Cells.Item(79, 8).Value = 55558364  # ALC!H79 (was 57695116)
$ws.Cells.Item(79, 9).Value = 65220120  # ALC!I79 (was 68184540)
$ws.Cells.Item(79, 11).Value = 65220120  # ALC!K79 (was 68184540)
$ws.Cells.Item(79, 13).Value = -65219028  # ALC!M79 (was -68183448)

$ws.Cells.Item(112, 8).Value = 1900  # ALC!H112 (was 2006.6666)
$ws.Cells.Item(112, 9).Value = 750  # ALC!I112 (was 766.6667)
$ws.Cells.Item(112, 10).Value = 2253.8462  # ALC!J112 (was 2316.6667)
$ws.Cells.Item(112, 11).Value = 2250  # ALC!K112 (was 2300.0001)
$ws.Cells.Item(112, 12).Value = 6761.5386  # ALC!L112 (was 6950.000100000001)
$ws.Cells.Item(112, 13).Value = -1142  # ALC!M112 (was -1192.0001)
$ws.Cells.Item(112, 14).Value = -8977.5386  # ALC!N112 (was -9166.000100000001)

$ws.Cells.Item(116, 8).Value = 4337.231  # ALC!H116 (was 4245.231)
$ws.Cells.Item(116, 9).Value = 3923.125  # ALC!I116 (was 4098.75)
$ws.Cells.Item(116, 10).Value = 4999.8  # ALC!J116 (was 4479.6)
$ws.Cells.Item(116, 11).Value = 3923.125  # ALC!K116 (was 4098.75)
$ws.Cells.Item(116, 12).Value = 4999.8  # ALC!L116 (was 4479.6)
$ws.Cells.Item(116, 13).Value = -481.125  # ALC!M116 (was -656.75)
$ws.Cells.Item(116, 14).Value = -11883.8  # ALC!N116 (was -11363.6)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3175.8  # ARM!H2 (was 2146.36)
$ws.Cells.Item(2, 9).Value = 2955.5  # ARM!I2 (was 1837.2354)
$ws.Cells.Item(2, 10).Value = 3427.5715  # ARM!J2 (was 2803.25)
$ws.Cells.Item(2, 11).Value = 2955.5  # ARM!K2 (was 1837.2354)
$ws.Cells.Item(2, 12).Value = 3427.5715  # ARM!L2 (was 2803.25)
$ws.Cells.Item(2, 13).Value = -2842.5  # ARM!M2 (was -1724.2354)
$ws.Cells.Item(2, 14).Value = -3653.5715  # ARM!N2 (was -3029.25)

$ws.Cells.Item(32, 8).Value = 3616.96  # ARM!H32 (was 4303.96)
$ws.Cells.Item(32, 9).Value = 3616.96  # ARM!I32 (was 4195.9194)
$ws.Cells.Item(32, 10).Value = 0  # ARM!J32 (was 15000)
$ws.Cells.Item(32, 11).Value = 3616.96  # ARM!K32 (was 4195.9194)
$ws.Cells.Item(32, 12).Value = 0  # ARM!L32 (was 15000)
$ws.Cells.Item(32, 13).Value = -3329.96  # ARM!M32 (was -3908.9194)
$ws.Cells.Item(32, 14).ClearContents()  # ARM!N32 (remove, was -15574)

$ws.Cells.Item(63, 8).Value = 2406.1875  # ARM!H63 (was 2212.1853)
$ws.Cells.Item(63, 9).Value = 2290.818  # ARM!I63 (was 2127.348)
$ws.Cells.Item(63, 10).Value = 2660  # ARM!J63 (was 2700)
$ws.Cells.Item(63, 11).Value = 2290.818  # ARM!K63 (was 2127.348)
$ws.Cells.Item(63, 12).Value = 2660  # ARM!L63 (was 2700)
$ws.Cells.Item(63, 13).Value = -1604.818  # ARM!M63 (was -1441.348)
$ws.Cells.Item(63, 14).Value = -4032  # ARM!N63 (was -4072)

$ws.Cells.Item(66, 8).Value = 2406.1875  # ARM!H66 (was 2212.1853)
$ws.Cells.Item(66, 9).Value = 2290.818  # ARM!I66 (was 2127.348)
$ws.Cells.Item(66, 10).Value = 2660  # ARM!J66 (was 2700)
$ws.Cells.Item(66, 11).Value = 11454.09  # ARM!K66 (was 10636.74)
$ws.Cells.Item(66, 12).Value = 13300  # ARM!L66 (was 13500)
$ws.Cells.Item(66, 13).Value = -8022.09  # ARM!M66 (was -7204.74)
$ws.Cells.Item(66, 14).Value = -20164  # ARM!N66 (was -20364)

$ws.Cells.Item(97, 8).Value = 3799.5881  # ARM!H97 (was 2108.5625)
$ws.Cells.Item(97, 9).Value = 5227.909  # ARM!I97 (was 2249.5557)
$ws.Cells.Item(97, 10).Value = 1181  # ARM!J97 (was 1347.2)
$ws.Cells.Item(97, 11).Value = 5227.909  # ARM!K97 (was 2249.5557)
$ws.Cells.Item(97, 12).Value = 1181  # ARM!L97 (was 1347.2)
$ws.Cells.Item(97, 13).Value = -4731.909  # ARM!M97 (was -1753.5557)
$ws.Cells.Item(97, 14).Value = -2173  # ARM!N97 (was -2339.2)

$ws.Cells.Item(116, 8).Value = 3175.8  # ARM!H116 (was 2146.36)
$ws.Cells.Item(116, 9).Value = 2955.5  # ARM!I116 (was 1837.2354)
$ws.Cells.Item(116, 10).Value = 3427.5715  # ARM!J116 (was 2803.25)
$ws.Cells.Item(116, 11).Value = 2955.5  # ARM!K116 (was 1837.2354)
$ws.Cells.Item(116, 12).Value = 3427.5715  # ARM!L116 (was 2803.25)
$ws.Cells.Item(116, 13).Value = -661.5  # ARM!M116 (was 456.7646)
$ws.Cells.Item(116, 14).Value = -8015.5715  # ARM!N116 (was -7391.25)

$ws.Cells.Item(122, 8).Value = 1343.3572  # ARM!H122 (was 1370.4615)
$ws.Cells.Item(122, 10).Value = 1455  # ARM!J122 (was 1513)
$ws.Cells.Item(122, 12).Value = 4365  # ARM!L122 (was 4539)
$ws.Cells.Item(122, 14).Value = -9265  # ARM!N122 (was -9439)

$ws.Cells.Item(132, 8).Value = 944.5  # ARM!H132 (was 957.28815)
$ws.Cells.Item(132, 9).Value = 807.05884  # ARM!I132 (was 819.4)
$ws.Cells.Item(132, 11).Value = 2421.17652  # ARM!K132 (was 2458.2)
$ws.Cells.Item(132, 13).Value = 108.82348  # ARM!M132 (was 71.80000000000018)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3175.8  # BSM!H3 (was 2146.36)
$ws.Cells.Item(3, 9).Value = 2955.5  # BSM!I3 (was 1837.2354)
$ws.Cells.Item(3, 10).Value = 3427.5715  # BSM!J3 (was 2803.25)
$ws.Cells.Item(3, 11).Value = 2955.5  # BSM!K3 (was 1837.2354)
$ws.Cells.Item(3, 12).Value = 3427.5715  # BSM!L3 (was 2803.25)
$ws.Cells.Item(3, 13).Value = -2841.5  # BSM!M3 (was -1723.2354)
$ws.Cells.Item(3, 14).Value = -3655.5715  # BSM!N3 (was -3031.25)

$ws.Cells.Item(7, 8).Value = 3750  # BSM!H7 (was 2033.3334)
$ws.Cells.Item(7, 9).Value = 3750  # BSM!I7 (was 2033.3334)
$ws.Cells.Item(7, 11).Value = 3750  # BSM!K7 (was 2033.3334)
$ws.Cells.Item(7, 13).Value = -3637  # BSM!M7 (was -1920.3334)

$ws.Cells.Item(10, 8).Value = 3000  # BSM!H10 (was 2495)
$ws.Cells.Item(10, 9).Value = 3000  # BSM!I10 (was 2495)
$ws.Cells.Item(10, 11).Value = 3000  # BSM!K10 (was 2495)
$ws.Cells.Item(10, 13).Value = -2860  # BSM!M10 (was -2355)

$ws.Cells.Item(17, 8).Value = 0  # BSM!H17 (was 3000)
$ws.Cells.Item(17, 10).Value = 0  # BSM!J17 (was 3000)
$ws.Cells.Item(17, 12).Value = 0  # BSM!L17 (was 3000)
$ws.Cells.Item(17, 14).ClearContents()  # BSM!N17 (remove, was -3344)

$ws.Cells.Item(99, 8).Value = 52633720  # BSM!H99 (was 19232274)
$ws.Cells.Item(99, 9).Value = 90911690  # BSM!I99 (was 33334822)
$ws.Cells.Item(99, 10).Value = 1512.5  # BSM!J99 (was 1527.2727)
$ws.Cells.Item(99, 11).Value = 90911690  # BSM!K99 (was 33334822)
$ws.Cells.Item(99, 12).Value = 1512.5  # BSM!L99 (was 1527.2727)
$ws.Cells.Item(99, 13).Value = -90910192  # BSM!M99 (was -33333324)
$ws.Cells.Item(99, 14).Value = -4508.5  # BSM!N99 (was -4523.2727)

$ws.Cells.Item(105, 8).Value = 4647.6763  # BSM!H105 (was 4607.1714)
$ws.Cells.Item(105, 9).Value = 4364.091  # BSM!I105 (was 4185.6)
$ws.Cells.Item(105, 10).Value = 5167.5835  # BSM!J105 (was 5661.1)
$ws.Cells.Item(105, 11).Value = 4364.091  # BSM!K105 (was 4185.6)
$ws.Cells.Item(105, 12).Value = 5167.5835  # BSM!L105 (was 5661.1)
$ws.Cells.Item(105, 13).Value = -2617.091  # BSM!M105 (was -2438.6)
$ws.Cells.Item(105, 14).Value = -8661.5835  # BSM!N105 (was -9155.1)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1073.3077  # CRP!H22 (was 1132.6666)
$ws.Cells.Item(22, 9).Value = 2168.6  # CRP!I22 (was 1860.3334)
$ws.Cells.Item(22, 10).Value = 388.75  # CRP!J22 (was 405)
$ws.Cells.Item(22, 11).Value = 2168.6  # CRP!K22 (was 1860.3334)
$ws.Cells.Item(22, 12).Value = 388.75  # CRP!L22 (was 405)
$ws.Cells.Item(22, 13).Value = -1818.6  # CRP!M22 (was -1510.3334)
$ws.Cells.Item(22, 14).Value = -1088.75  # CRP!N22 (was -1105)

$ws.Cells.Item(31, 8).Value = 39041.15  # CRP!H31 (was 33328.28)
$ws.Cells.Item(31, 9).Value = 1820.55  # CRP!I31 (was 1845.25)
$ws.Cells.Item(31, 10).Value = 145385.72  # CRP!J31 (was 85800)
$ws.Cells.Item(31, 11).Value = 1820.55  # CRP!K31 (was 1845.25)
$ws.Cells.Item(31, 12).Value = 145385.72  # CRP!L31 (was 85800)
$ws.Cells.Item(31, 13).Value = -1525.55  # CRP!M31 (was -1550.25)
$ws.Cells.Item(31, 14).Value = -145975.72  # CRP!N31 (was -86390)

$ws.Cells.Item(34, 8).Value = 39041.15  # CRP!H34 (was 33328.28)
$ws.Cells.Item(34, 9).Value = 1820.55  # CRP!I34 (was 1845.25)
$ws.Cells.Item(34, 10).Value = 145385.72  # CRP!J34 (was 85800)
$ws.Cells.Item(34, 11).Value = 1820.55  # CRP!K34 (was 1845.25)
$ws.Cells.Item(34, 12).Value = 145385.72  # CRP!L34 (was 85800)
$ws.Cells.Item(34, 13).Value = -1618.55  # CRP!M34 (was -1643.25)
$ws.Cells.Item(34, 14).Value = -145789.72  # CRP!N34 (was -86204)

$ws.Cells.Item(58, 8).Value = 2547.3708  # CRP!H58 (was 2289.845)
$ws.Cells.Item(58, 9).Value = 911.1429  # CRP!I58 (was 822.93335)
$ws.Cells.Item(58, 10).Value = 4668.407  # CRP!J58 (was 4828.731)
$ws.Cells.Item(58, 11).Value = 911.1429  # CRP!K58 (was 822.93335)
$ws.Cells.Item(58, 12).Value = 4668.407  # CRP!L58 (was 4828.731)
$ws.Cells.Item(58, 13).Value = -708.1429  # CRP!M58 (was -619.93335)
$ws.Cells.Item(58, 14).Value = -5074.407  # CRP!N58 (was -5234.731)

$ws.Cells.Item(62, 8).Value = 4267.1113  # CRP!H62 (was 4857)
$ws.Cells.Item(62, 9).Value = 2476  # CRP!I62 (was 2749.5)
$ws.Cells.Item(62, 11).Value = 2476  # CRP!K62 (was 2749.5)
$ws.Cells.Item(62, 13).Value = -1852  # CRP!M62 (was -2125.5)

$ws.Cells.Item(65, 8).Value = 4267.1113  # CRP!H65 (was 4857)
$ws.Cells.Item(65, 9).Value = 2476  # CRP!I65 (was 2749.5)
$ws.Cells.Item(65, 11).Value = 12380  # CRP!K65 (was 13747.5)
$ws.Cells.Item(65, 13).Value = -9260  # CRP!M65 (was -10627.5)

$ws.Cells.Item(132, 8).Value = 1472.0312  # CRP!H132 (was 1356.6945)
$ws.Cells.Item(132, 9).Value = 960.5263  # CRP!I132 (was 904.0952)
$ws.Cells.Item(132, 10).Value = 2219.6155  # CRP!J132 (was 1990.3334)
$ws.Cells.Item(132, 11).Value = 2881.5789  # CRP!K132 (was 2712.2856)
$ws.Cells.Item(132, 12).Value = 6658.8465  # CRP!L132 (was 5971.0002)
$ws.Cells.Item(132, 13).Value = -351.5789  # CRP!M132 (was -182.2856000000002)
$ws.Cells.Item(132, 14).Value = -11718.8465  # CRP!N132 (was -11031.0002)

$ws.Cells.Item(134, 8).Value = 20001144  # CRP!H134 (was 10000741)
$ws.Cells.Item(134, 9).Value = 1117.4546  # CRP!I134 (was 729.75)
$ws.Cells.Item(134, 10).Value = 166668000  # CRP!J134 (was 250001000)
$ws.Cells.Item(134, 11).Value = 3352.3638  # CRP!K134 (was 2189.25)
$ws.Cells.Item(134, 12).Value = 500004000  # CRP!L134 (was 750003000)
$ws.Cells.Item(134, 13).Value = -817.3638000000001  # CRP!M134 (was 345.75)
$ws.Cells.Item(134, 14).Value = -500009070  # CRP!N134 (was -750008070)

$ws.Cells.Item(136, 8).Value = 2547.3708  # CRP!H136 (was 2289.845)
$ws.Cells.Item(136, 9).Value = 911.1429  # CRP!I136 (was 822.93335)
$ws.Cells.Item(136, 10).Value = 4668.407  # CRP!J136 (was 4828.731)
$ws.Cells.Item(136, 11).Value = 2733.4287  # CRP!K136 (was 2468.80005)
$ws.Cells.Item(136, 12).Value = 14005.221  # CRP!L136 (was 14486.193)
$ws.Cells.Item(136, 13).Value = -183.4287000000004  # CRP!M136 (was 81.19995000000017)
$ws.Cells.Item(136, 14).Value = -19105.221  # CRP!N136 (was -19586.193)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4475.1055  # GSM!H70 (was 4300.6816)
$ws.Cells.Item(70, 9).Value = 4014.2856  # GSM!I70 (was 3788.6667)
$ws.Cells.Item(70, 10).Value = 4743.9165  # GSM!J70 (was 4655.154)
$ws.Cells.Item(70, 11).Value = 4014.2856  # GSM!K70 (was 3788.6667)
$ws.Cells.Item(70, 12).Value = 4743.9165  # GSM!L70 (was 4655.154)
$ws.Cells.Item(70, 13).Value = -3744.2856  # GSM!M70 (was -3518.6667)
$ws.Cells.Item(70, 14).Value = -5283.9165  # GSM!N70 (was -5195.154)

$ws.Cells.Item(73, 8).Value = 4475.1055  # GSM!H73 (was 4300.6816)
$ws.Cells.Item(73, 9).Value = 4014.2856  # GSM!I73 (was 3788.6667)
$ws.Cells.Item(73, 10).Value = 4743.9165  # GSM!J73 (was 4655.154)
$ws.Cells.Item(73, 11).Value = 4014.2856  # GSM!K73 (was 3788.6667)
$ws.Cells.Item(73, 12).Value = 4743.9165  # GSM!L73 (was 4655.154)
$ws.Cells.Item(73, 13).Value = -3078.2856  # GSM!M73 (was -2852.6667)
$ws.Cells.Item(73, 14).Value = -6615.9165  # GSM!N73 (was -6527.154)

$ws.Cells.Item(80, 8).Value = 3507.7273  # GSM!H80 (was 3239.6875)
$ws.Cells.Item(80, 9).Value = 3508.5  # GSM!I80 (was 3233.4614)
$ws.Cells.Item(80, 10).Value = 3500  # GSM!J80 (was 3266.6667)
$ws.Cells.Item(80, 11).Value = 3508.5  # GSM!K80 (was 3233.4614)
$ws.Cells.Item(80, 12).Value = 3500  # GSM!L80 (was 3266.6667)
$ws.Cells.Item(80, 13).Value = -2510.5  # GSM!M80 (was -2235.4614)
$ws.Cells.Item(80, 14).Value = -5496  # GSM!N80 (was -5262.6667)

$ws.Cells.Item(83, 8).Value = 3507.7273  # GSM!H83 (was 3239.6875)
$ws.Cells.Item(83, 9).Value = 3508.5  # GSM!I83 (was 3233.4614)
$ws.Cells.Item(83, 10).Value = 3500  # GSM!J83 (was 3266.6667)
$ws.Cells.Item(83, 11).Value = 17542.5  # GSM!K83 (was 16167.307)
$ws.Cells.Item(83, 12).Value = 17500  # GSM!L83 (was 16333.3335)
$ws.Cells.Item(83, 13).Value = -12550.5  # GSM!M83 (was -11175.307)
$ws.Cells.Item(83, 14).Value = -27484  # GSM!N83 (was -26317.3335)

$ws.Cells.Item(107, 8).Value = 458.85715  # GSM!H107 (was 361.6842)
$ws.Cells.Item(107, 9).Value = 283.875  # GSM!I107 (was 256.77777)
$ws.Cells.Item(107, 10).Value = 692.1667  # GSM!J107 (was 456.1)
$ws.Cells.Item(107, 11).Value = 283.875  # GSM!K107 (was 256.77777)
$ws.Cells.Item(107, 12).Value = 692.1667  # GSM!L107 (was 456.1)
$ws.Cells.Item(107, 13).Value = 1636.125  # GSM!M107 (was 1663.22223)
$ws.Cells.Item(107, 14).Value = -4532.1667  # GSM!N107 (was -4296.1)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 8571.467  # LTW!H122 (was 10225.417)
$ws.Cells.Item(122, 9).Value = 12796.333  # LTW!I122 (was 15928.571)
$ws.Cells.Item(122, 10).Value = 2234.1667  # LTW!J122 (was 2241)
$ws.Cells.Item(122, 11).Value = 38388.999  # LTW!K122 (was 47785.713)
$ws.Cells.Item(122, 12).Value = 6702.500100000001  # LTW!L122 (was 6723)
$ws.Cells.Item(122, 13).Value = -35938.999  # LTW!M122 (was -45335.713)
$ws.Cells.Item(122, 14).Value = -11602.5001  # LTW!N122 (was -11623)

$ws.Cells.Item(136, 8).Value = 3496.362  # LTW!H136 (was 3855.255)
$ws.Cells.Item(136, 9).Value = 2039.6945  # LTW!I136 (was 2319.2415)
$ws.Cells.Item(136, 11).Value = 6119.083500000001  # LTW!K136 (was 6957.7245)
$ws.Cells.Item(136, 13).Value = -3569.083500000001  # LTW!M136 (was -4407.7245)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 29666.666  # WVR!H16 (was 0)
$ws.Cells.Item(16, 10).Value = 29666.666  # WVR!J16 (was 0)
$ws.Cells.Item(16, 12).Value = 29666.666  # WVR!L16 (was 0)
$ws.Cells.Item(16, 14).Value = -30250.666  # WVR!N16 (was None)

$ws.Cells.Item(132, 8).Value = 1066.2927  # WVR!H132 (was 1073.7561)
$ws.Cells.Item(132, 9).Value = 967  # WVR!I132 (was 977.9286)
$ws.Cells.Item(132, 11).Value = 2901  # WVR!K132 (was 2933.7858)
$ws.Cells.Item(132, 13).Value = -371  # WVR!M132 (was -403.7857999999997)

